$wb = $excel.ActiveWorkbook

# --- Sheet references -------------------------------------------------
$wsRevHistory  = $wb.Worksheets.Item("Revision History")
$wsInstructions = $wb.Worksheets.Item("Instructions for Use")
$wsChangeReq   = $wb.Worksheets.Item("QMS_X.X Change Requirements")

# --- Text / content updates --------------------------------------------
# Revision History: bump the template version string
$wsRevHistory.Range("A3").Value = "Template Version -3"

# QMS Change Requirements: rename "SVN Repository Rev." -> "Repository Rev."
# (this updates both the table header cell and the table column definition)
$wsChangeReq.Range("H5").Value = "Repository Rev."

$tbl = $wsChangeReq.ListObjects.Item(1)
$tbl.ListColumns.Item(8).Name = "Repository Rev."

# Swap the other note back to "Points to cover in Post Release training"
$wsChangeReq.Range("A14").Value = "Points to cover in Post Release training"

# --- View / selection state ---------------------------------------------
# Visit the change-requirements sheet and leave the selection on G13
# (also clears its previous tabSelected/topLeftCell scroll state)
$wsChangeReq.Activate() | Out-Null
$wsChangeReq.Range("G13").Select() | Out-Null

# Finish on "Instructions for Use" with A3 selected - this becomes the
# active / selected tab when the workbook is saved
$wsInstructions.Activate() | Out-Null
$wsInstructions.Range("A3").Select() | Out-Null
